$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Update Version value
$ws.Range("B3").Value = "6.0.0"

# Update Date value
$ws.Range("B8").Value = "2022-01-21T20:46:54+00:00"

# Update Publisher value
$ws.Range("B9").Value = "Alvearie Team"

# Remove the second duplicate "Contact" row (row 11) first, so the
# remaining "Contact" row (row 10) can be turned into "Jurisdiction"
$ws.Rows.Item(11).Delete()

# Update row 10 ("Contact" -> "Jurisdiction")
$ws.Range("A10").Value = "Jurisdiction"
$ws.Range("B10").Value = "United States of America"
